# Swap the data of row 19 and row 20 (the two observation records got
# reordered). Only the cells whose content actually differs between the
# two rows are touched; cells already identical between row 19 and row 20
# are left alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 19: becomes what row 20 used to contain ----
$ws.Range("A19").Value = 111768503
$ws.Range("B19").Value = 88966
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 5754
$ws.Range("F19").Value = "Gultoppig fingersvamp"
$ws.Range("G19").Value = "Ramaria testaceoflava"
$ws.Range("H19").Value = "(Bres.) Corner"
$ws.Range("I19").NumberFormat = "@"   # keep "Antal" as text, matching the source data
$ws.Range("I19").Value = "20"
$ws.Range("J19").Value = "fruktkroppar"
$ws.Range("L19").Value = ""        # row 20 had no L cell -> remove it from row 19
$ws.Range("Q19").Value = 525545.3455456314
$ws.Range("R19").Value = 6727837.787189188
$ws.Range("Z19").Value = "15:22"
$ws.Range("AB19").Value = "15:22"

# ---- Row 20: becomes what row 19 used to contain ----
$ws.Range("A20").Value = 111768476
$ws.Range("B20").Value = 96348
$ws.Range("D20").Value = "VU"
$ws.Range("E20").Value = 220787
$ws.Range("F20").Value = "Knärot"
$ws.Range("G20").Value = "Goodyera repens"
$ws.Range("H20").Value = "(L.) R. Br."
$ws.Range("I20").NumberFormat = "@"   # keep "Antal" as text, matching the source data
$ws.Range("I20").Value = "25"
$ws.Range("J20").Value = "plantor/tuvor"
$ws.Range("L20").Value = ""        # row 19 had an (empty) L cell -> add it to row 20
$ws.Range("Q20").Value = 525546.5036804043
$ws.Range("R20").Value = 6727881.884716956
$ws.Range("Z20").Value = "15:21"
$ws.Range("AB20").Value = "15:21"
